$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Static offsets: zero out every existing station's offset (column C, rows 1-30) ---
$ws.Range("C1:C30").Value = 0

# --- 2. New stations appended as rows 31-33 ---

# Row 31: USACE station 01480 (leading-zero code -> must be stored as text)
$ws.Range("A31").Value = "USACE"
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "01480"
$ws.Range("C31").NumberFormat = "0.00"
$ws.Range("C31").Value = 0

# Row 32: USACE station 76560 (plain numeric code)
$ws.Range("A32").Value = "USACE"
$ws.Range("B32").Value = 76560
$ws.Range("C32").NumberFormat = "0.00"
$ws.Range("C32").Value = 0

# Row 33: USGS station 073814675 (leading-zero code -> must be stored as text)
$ws.Range("A33").Value = "USGS"
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "073814675"
$ws.Range("C33").NumberFormat = "0.00"
$ws.Range("C33").Value = 0

# --- 3. Existing USACE station 82740 (row 7) replaced by 82742 ---
# Done after the new rows are written so the shared-string table matches
# the order strings were actually introduced.
$ws.Range("B7").Value = "82742"

# --- 4. Update the view: select the whole (now taller) offset column ---
$ws.Range("A1").Select() | Out-Null
$ws.Range("C1:C33").Select() | Out-Null
